$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Sir, I’ll leave it to you to find the culprit."
$ws.Range("B5").Value = "When was the last time you saw the Lord?"
$ws.Range("B6").Value = "Sometime after 5 PM. I was walking down the corridor and happened to see the master coming out of the main hall."
$ws.Range("B8").Value = "He often guided me in martial arts before."
$ws.Range("B12").Value = " <color=#00CC00>(Butler He and Ling both mentioned this in their statements.)</color>"
$ws.Range("B14").Value = " <color=#00CC00>(I recall they did indeed arrive at the banquet hall at the same time.)</color>"
$ws.Range("B16").Value = "Does going to the restroom in the banquet hall count? I was gone for about 15 min."

$ws.Rows.Item(5).RowHeight = 17

$ws.Range("B22").Select()
